$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old "me" / "mememe" values so they get garbage-collected
# from the shared string table before the new strings are introduced.
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""

# Introduce the new shared strings in the order that reproduces the
# target shared-string table ordering (VQSR, PASS, FAIL, then the
# duplicate "PASS" cells that reuse the same entry).
$ws.Range("E5").Value = "VQSR"
$ws.Range("E2").Value = "PASS"
$ws.Range("E4").Value = "FAIL"
$ws.Range("E3").Value = "PASS"
$ws.Range("E6").Value = "PASS"

# Apply the new compare-method font (black RGB Calibri) to the new
# "compare method" / result columns, matching the new cellXfs entry.
$ws.Range("E2:F6").Font.Color = 0

# Fill in the new numeric compare values.
$ws.Range("F2").Value = 0.2
$ws.Range("F3").Value = 0.01
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = 0.04
$ws.Range("F6").Value = 0.0004

# Update the active cell selection on the sheet.
$ws.Range("G14").Select()
